$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers (row 3): Elemento 1..6 -> Fastball, Riseball, Dropball, Curveball, ChangeUp, Screwball
$ws.Range("B3").Value = "Fastball"
$ws.Range("C3").Value = "Riseball"
$ws.Range("D3").Value = "Dropball"
$ws.Range("E3").Value = "Curveball"
$ws.Range("F3").Value = "ChangeUp"
$ws.Range("G3").Value = "Screwball"

# Row headers (column A, rows 4-9): Caracteristica 1..6 -> Velocidad, Posicion Cuerpo, Traza, Direccion, Efecto, Agarre
$ws.Range("A4").Value = "Velocidad"
$ws.Range("A5").Value = "Posición Cuerpo"
$ws.Range("A6").Value = "Traza"
$ws.Range("A7").Value = "Dirección"
$ws.Range("A8").Value = "Efecto"
$ws.Range("A9").Value = "Agarre"

# Selected cell moves from A17 to A10
$ws.Range("A10").Select()
